$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '29.215.64'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +1.15%  '

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.938.99'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +2.40%  '

# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.004'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '326.38'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.08%  '

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.31%  '

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4628'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.73%  '

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3903'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.08%  '

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.07876'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +0.46%  '

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.9995'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +0.93%  '

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '22.25'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +1.37%  '

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.923.86'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -1.05%  '

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '5.821'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +2.17%  '

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '7.094'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +0.69%  '

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.07060'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +1.51%  '

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '88.06'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +0.06%  '

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.13%  '

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.000009959'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -0.07%  '

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '17.20'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +1.57%  '

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.14%  '

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '29.234.01'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +1.11%  '

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.478'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +3.33%  '

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '11.22'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +2.24%  '

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.159.89'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.08%  '

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.100'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +1.81%  '

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '156.31'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.09%  '

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '19.53'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +1.36%  '

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '5.904'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -0.72%  '

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '119.07'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +1.28%  '

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.885'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -2.77%  '

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.09346'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.03%  '

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.8948'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -2.14%  '

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '5.231'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -1.27%  '

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.326'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -0.80%  '

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '3.134'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -4.60%  '

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.05782'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +0.14%  '

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.173'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -1.56%  '

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.02102'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +1.25%  '

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '1.0000'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.25%  '

# Row 40
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '7.701'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -1.21%  '

# Row 41
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.5719'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +0.58%  '

# Row 42
$ws.Range("B42").Value = 'PEPE'
$ws.Range("C42").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.000003141'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +97.87%  '

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.1819'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +2.44%  '

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '9.768'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -0.26%  '

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '11.98'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.73%  '

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '2.217'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -1.84%  '

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.5358'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +0.01%  '

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.06950'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -1.31%  '

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.855'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +0.83%  '

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '2.589'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +2.49%  '

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '113.21'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +0.65%  '

